$d = $word.ActiveDocument

# Work only within the document's first paragraph so the later historical
# "...Observed till July 2022" / "...31st Oct 2021" entries are untouched.
$para = $d.Paragraphs(1)
$rng = $para.Range

# Step 1: merge the first three runs ("TS Krama Paatam – TS 1.6 " + "Tamil" +
# " Corrections – Observed till ") into the wording used by the new title,
# leaving the trailing placeholder run alone for now.
$rng.Find.Execute("TS Krama Paatam – TS 1.6 Tamil Corrections – Observed till ", $true, $false, $false, $false, $false, $true, 1, $false, "TS Krama Paatam – TS 1.6 Tamil Corrections – Observed till ", 2)

# Step 2: replace the red-highlighted "??????" placeholder with the new date.
$rng2 = $para.Range
$rng2.Find.Execute("??????", $true, $false, $false, $false, $false, $true, 1, $false, "31st Jan 2026", 2)

# Step 3: clear the red highlight left over on the placeholder's old run so
# the date run matches the rest of the heading (no highlight formatting).
$rng3 = $para.Range
$rng3.Find.ClearFormatting()
$rng3.Find.Text = "31st Jan 2026"
$rng3.Find.Execute()
if ($rng3.Find.Found) {
    $rng3.HighlightColorIndex = 0
}
